$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3164.889
$ws.Range("I40").Value = 1296.6666
$ws.Range("J40").Value = 4099
$ws.Range("K40").Value = 1296.6666
$ws.Range("L40").Value = 4099
$ws.Range("M40").Value = -1121.6666
$ws.Range("N40").Value = -4449
$ws.Range("H69").Value = 71431976
$ws.Range("J69").Value = 71431976
$ws.Range("L69").Value = 214295928
$ws.Range("N69").Value = -214297676
$ws.Range("H72").Value = 71431976
$ws.Range("J72").Value = 71431976
$ws.Range("L72").Value = 642887784
$ws.Range("N72").Value = -642896520
$ws.Range("H101").Value = 2338
$ws.Range("I101").Value = 1073.2
$ws.Range("J101").Value = 5500
$ws.Range("K101").Value = 3219.6
$ws.Range("L101").Value = 16500
$ws.Range("M101").Value = -1597.6
$ws.Range("N101").Value = -19744
$ws.Range("H112").Value = 5234.075
$ws.Range("J112").Value = 5618.5947
$ws.Range("L112").Value = 16855.7841
$ws.Range("N112").Value = -19071.7841
$ws.Range("H116").Value = 7411.522
$ws.Range("I116").Value = 10383.462
$ws.Range("J116").Value = 3548
$ws.Range("K116").Value = 10383.462
$ws.Range("L116").Value = 3548
$ws.Range("M116").Value = -6941.462
$ws.Range("N116").Value = -10432
$ws.Range("H134").Value = 50186.668
$ws.Range("J134").Value = 50186.668
$ws.Range("L134").Value = 50186.668
$ws.Range("N134").Value = -60326.668
$ws.Range("H138").Value = 3730.5544
$ws.Range("I138").Value = 6084.25
$ws.Range("J138").Value = 3506.3928
$ws.Range("K138").Value = 18252.75
$ws.Range("L138").Value = 10519.1784
$ws.Range("M138").Value = -13112.75
$ws.Range("N138").Value = -20799.1784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 23018.625
$ws.Range("J37").Value = 25735.572
$ws.Range("L37").Value = 25735.572
$ws.Range("N37").Value = -26281.572
$ws.Range("H44").Value = 26749.25
$ws.Range("H55").Value = 26999.2
$ws.Range("H132").Value = 1401683.8
$ws.Range("I132").Value = 2507.457
$ws.Range("J132").Value = 3850242.5
$ws.Range("K132").Value = 7522.370999999999
$ws.Range("L132").Value = 11550727.5
$ws.Range("M132").Value = -4992.370999999999
$ws.Range("N132").Value = -11555787.5
$ws.Range("H138").Value = 138457.2
$ws.Range("J138").Value = 138457.2
$ws.Range("L138").Value = 138457.2
$ws.Range("N138").Value = -148737.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 22728694
$ws.Range("I20").Value = 1316.6207
$ws.Range("J20").Value = 66668290
$ws.Range("K20").Value = 1316.6207
$ws.Range("L20").Value = 66668290
$ws.Range("M20").Value = -1069.6207
$ws.Range("N20").Value = -66668784
$ws.Range("H99").Value = 1716.5
$ws.Range("I99").Value = 1270
$ws.Range("J99").Value = 1984.4
$ws.Range("K99").Value = 1270
$ws.Range("L99").Value = 1984.4
$ws.Range("M99").Value = 228
$ws.Range("N99").Value = -4980.4
$ws.Range("H134").Value = 2636
$ws.Range("I134").Value = 2652.1428
$ws.Range("J134").Value = 2598.3333
$ws.Range("K134").Value = 7956.428400000001
$ws.Range("L134").Value = 7794.999899999999
$ws.Range("M134").Value = -5421.428400000001
$ws.Range("N134").Value = -12864.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1200
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 1000
$ws.Range("M8").Value = -860
$ws.Range("H22").Value = 452.8
$ws.Range("I22").Value = 291.66666
$ws.Range("J22").Value = 560.2222
$ws.Range("K22").Value = 291.66666
$ws.Range("L22").Value = 560.2222
$ws.Range("M22").Value = 58.33334000000002
$ws.Range("N22").Value = -1260.2222
$ws.Range("H31").Value = 5521.0405
$ws.Range("I31").Value = 1733.8
$ws.Range("J31").Value = 7453.306
$ws.Range("K31").Value = 1733.8
$ws.Range("L31").Value = 7453.306
$ws.Range("M31").Value = -1438.8
$ws.Range("N31").Value = -8043.306
$ws.Range("H34").Value = 5521.0405
$ws.Range("I34").Value = 1733.8
$ws.Range("J34").Value = 7453.306
$ws.Range("K34").Value = 1733.8
$ws.Range("L34").Value = 7453.306
$ws.Range("M34").Value = -1531.8
$ws.Range("N34").Value = -7857.306
$ws.Range("H122").Value = 2092.258
$ws.Range("I122").Value = 1834.7368
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5504.2104
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3054.2104
$ws.Range("N122").Value = -12400
$ws.Range("H123").Value = 36926.668
$ws.Range("J123").Value = 36926.668
$ws.Range("L123").Value = 36926.668
$ws.Range("N123").Value = -46726.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 8250
$ws.Range("I120").Value = 2500
$ws.Range("J120").Value = 14000
$ws.Range("K120").Value = 7500
$ws.Range("L120").Value = 42000
$ws.Range("M120").Value = -2662
$ws.Range("N120").Value = -51676
$ws.Range("H131").Value = 3366.5833
$ws.Range("I131").Value = 477.14285
$ws.Range("J131").Value = 4556.353
$ws.Range("K131").Value = 1431.42855
$ws.Range("L131").Value = 13669.059
$ws.Range("M131").Value = 3608.57145
$ws.Range("N131").Value = -23749.059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 11100170
$ws.Range("I14").Value = 11100170
$ws.Range("K14").Value = 11100170
$ws.Range("M14").Value = -11100002
$ws.Range("H34").Value = 32000
$ws.Range("J34").Value = 32000
$ws.Range("L34").Value = 32000
$ws.Range("N34").Value = -32536
$ws.Range("H70").Value = 5706.353
$ws.Range("I70").Value = 5646.769
$ws.Range("J70").Value = 5900
$ws.Range("K70").Value = 5646.769
$ws.Range("L70").Value = 5900
$ws.Range("M70").Value = -5376.769
$ws.Range("N70").Value = -6440
$ws.Range("H73").Value = 5706.353
$ws.Range("I73").Value = 5646.769
$ws.Range("J73").Value = 5900
$ws.Range("K73").Value = 5646.769
$ws.Range("L73").Value = 5900
$ws.Range("M73").Value = -4710.769
$ws.Range("N73").Value = -7772
$ws.Range("H76").Value = 32000
$ws.Range("J76").Value = 32000
$ws.Range("L76").Value = 32000
$ws.Range("N76").Value = -32630
$ws.Range("H79").Value = 32000
$ws.Range("J79").Value = 32000
$ws.Range("L79").Value = 32000
$ws.Range("N79").Value = -34184
$ws.Range("H80").Value = 1881241
$ws.Range("I80").Value = 3001068.2
$ws.Range("J80").Value = 201500
$ws.Range("K80").Value = 3001068.2
$ws.Range("L80").Value = 201500
$ws.Range("M80").Value = -3000070.2
$ws.Range("N80").Value = -203496
$ws.Range("H83").Value = 1881241
$ws.Range("I83").Value = 3001068.2
$ws.Range("J83").Value = 201500
$ws.Range("K83").Value = 15005341
$ws.Range("L83").Value = 1007500
$ws.Range("M83").Value = -15000349
$ws.Range("N83").Value = -1017484
$ws.Range("H102").Value = 1492.1852
$ws.Range("I102").Value = 1561.55
$ws.Range("J102").Value = 1294
$ws.Range("K102").Value = 1561.55
$ws.Range("L102").Value = 1294
$ws.Range("M102").Value = 60.45000000000005
$ws.Range("N102").Value = -4538
$ws.Range("H132").Value = 32264132
$ws.Range("I132").Value = 52639044
$ws.Range("J132").Value = 3853.0833
$ws.Range("K132").Value = 157917132
$ws.Range("L132").Value = 11559.2499
$ws.Range("M132").Value = -157914602
$ws.Range("N132").Value = -16619.2499
$ws.Range("H141").Value = 79395.664
$ws.Range("J141").Value = 79395.664
$ws.Range("L141").Value = 79395.664
$ws.Range("N141").Value = -89755.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2500
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 2500
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 2500
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -2948
$ws.Range("H82").Value = 2897
$ws.Range("I82").Value = 3052.8572
$ws.Range("J82").Value = 2533.3333
$ws.Range("K82").Value = 3052.8572
$ws.Range("L82").Value = 2533.3333
$ws.Range("M82").Value = -2691.8572
$ws.Range("N82").Value = -3255.3333
$ws.Range("H85").Value = 2897
$ws.Range("I85").Value = 3052.8572
$ws.Range("J85").Value = 2533.3333
$ws.Range("K85").Value = 3052.8572
$ws.Range("L85").Value = 2533.3333
$ws.Range("M85").Value = -1804.8572
$ws.Range("N85").Value = -5029.3333
$ws.Range("H132").Value = 3106.4
$ws.Range("I132").Value = 2621.889
$ws.Range("J132").Value = 3833.1667
$ws.Range("K132").Value = 7865.667
$ws.Range("L132").Value = 11499.5001
$ws.Range("M132").Value = -5335.667
$ws.Range("N132").Value = -16559.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2404.6875
$ws.Range("I122").Value = 2459.6155
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 7378.8465
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -4928.8465
$ws.Range("N122").Value = -11400.0001
